$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6545652718822623
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 3993.344853322108
$ws.Range("E2").Value = 14773364.14517103
$ws.Range("G2").Value = 14777358.44948087
